$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: update title and link
$ws.Range("D23").Value = "파이참Pycharm으로 원격 서버 SSH 연결해서 코드 수정하기 (윈도우)"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2959"

# Row 26: update title
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

# Row 36: fix capitalization in title (MAP -> Map)
$ws.Range("D36").Value = "History of Class Activation Map (CAM)"
